$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 75
$ws.Range("F3").Value = 200
$ws.Range("F4").Value = 68
$ws.Range("F5").Value = 1653
$ws.Range("F6").Value = 3253
$ws.Range("F7").Value = 809
$ws.Range("F8").Value = 2047
$ws.Range("F9").Value = 1961
$ws.Range("F10").Value = 1001
$ws.Range("F11").Value = 352
$ws.Range("F13").Value = 1608
$ws.Range("F14").Value = 346
$ws.Range("F16").Value = 12
$ws.Range("F17").Value = 72
$ws.Range("F18").Value = 59
$ws.Range("F19").Value = 1437
$ws.Range("F20").Value = 522
$ws.Range("F21").Value = 624
$ws.Range("F22").Value = 319
$ws.Range("F23").Value = 10679
$ws.Range("F24").Value = 9830
$ws.Range("F25").Value = 846
$ws.Range("F26").Value = 655
$ws.Range("F27").Value = 1834
$ws.Range("F28").Value = 146
$ws.Range("F29").Value = 427

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 36

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 61

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 75
$ws.Range("F3").Value = 61
$ws.Range("F4").Value = 200
$ws.Range("F6").Value = 68
$ws.Range("F7").Value = 1653
$ws.Range("F8").Value = 3253
$ws.Range("F9").Value = 809
$ws.Range("F10").Value = 2047
$ws.Range("F11").Value = 1961
$ws.Range("F12").Value = 1001
$ws.Range("F13").Value = 352
$ws.Range("F15").Value = 1608
$ws.Range("F16").Value = 346
$ws.Range("F18").Value = 12
$ws.Range("F20").Value = 72
$ws.Range("F21").Value = 36
$ws.Range("F22").Value = 59
$ws.Range("F23").Value = 1437
$ws.Range("F24").Value = 522
$ws.Range("F25").Value = 624
$ws.Range("F26").Value = 319
$ws.Range("F27").Value = 10679
$ws.Range("F28").Value = 9830
$ws.Range("F29").Value = 846
$ws.Range("F30").Value = 655
$ws.Range("F31").Value = 1834
$ws.Range("F34").Value = 146
$ws.Range("F35").Value = 427
